$d = $word.ActiveDocument

$replacements = @(
    @("2025-05-19 Monday", "2025-05-20 Tuesday"),
    @("58×26=1508", "60×17=1020"),
    @("81×55=4455", "89×94=8366"),
    @("68×46=3128", "14×76=1064"),
    @("32×85=2720", "90×88=7920"),
    @("65×80=5200", "59×21=1239"),
    @("75×77=5775", "60×26=1560"),
    @("79×63=4977", "23×93=2139"),
    @("92×35=3220", "87×75=6525"),
    @("54×89=4806", "66×32=2112"),
    @("25×21=525",  "70×97=6790"),
    @("93×84=7812", "52×82=4264"),
    @("73×69=5037", "28×55=1540"),
    @("99×22=2178", "47×70=3290"),
    @("48×97=4656", "62×66=4092"),
    @("61×62=3782", "90×12=1080"),
    @("26×82=2132", "59×79=4661"),
    @("55×40=2200", "23×93=2139"),
    @("59×56=3304", "98×37=3626"),
    @("62×23=1426", "77×81=6237"),
    @("82×17=1394", "47×87=4089"),
    @("93×47=4371", "83×33=2739"),
    @("88×20=1760", "45×94=4230"),
    @("66×77=5082", "88×34=2992"),
    @("83×60=4980", "11×11=121"),
    @("87×88=7656", "60×11=660")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
